$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new "Номер" column (A2:A5) the same bold/centered/bordered
# format already used by the header row (B1:D1) by copying B1's format.
$ws.Range("B1").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Test-run result rows: row number, login, result
$data = @(
    @(0, "FDTK1CB7140", "Неверный"),
    @(1, "DPTWNT8K140", "ПройденоУспех"),
    @(2, "фыв", "Неверный"),
    @(3, "asf", "Неверный")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2

    $ws.Range("A$row").Value = $data[$i][0]

    # Column B stays blank for every data row; force Excel to keep a real
    # (empty) text cell instead of silently clearing it by writing a
    # quote-prefixed empty value, then restore the default "Normal" style.
    $ws.Range("B$row").Value = "'"
    $ws.Range("B$row").Style = "Normal"

    $ws.Range("C$row").Value = $data[$i][1]
    $ws.Range("D$row").Value = $data[$i][2]
}
